$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at GT (shifts old GT->GU, old GU->GV),
# pushing "nom"/"url_produit" columns one to the right.
$ws.Columns("GT:GT").Insert()

# New header cell: latest snapshot timestamp.
$ws.Range("GT1").Value = "2026-02-06 09:32:56"

# Populate the new snapshot column (GT) with the same value as the
# previous last snapshot column (GS) for every data row, mirroring the
# "carry forward last known price" behaviour of the scraper.
$lastRow = 210
for ($r = 2; $r -le $lastRow; $r++) {
    $prevValue = $ws.Cells.Item($r, 201).Value2  # column GS
    if ($prevValue -ne "") {
        $ws.Cells.Item($r, 202).Value2 = $prevValue  # column GT (newly inserted)
    }
}
